$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column A (Id), Q (Ost), R (Nord) for rows 17-21.
# This reflects a re-sequencing / shift of the underlying source rows.
$updates = @{
    17 = @{ A = 111821928; Q = 550825.9503372401; R = 6681726.144349095 }
    18 = @{ A = 111821927; Q = 550819.8901872271; R = 6681733.007140613 }
    19 = @{ A = 111821923; Q = 550701.1291094749; R = 6681909.496304798 }
    20 = @{ A = 111821926; Q = 550846.2444635418; R = 6681625.195240833 }
    21 = @{ A = 111821924; Q = 550675.3931295178; R = 6681937.422269406 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("A$row").Value = $vals.A
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("R$row").Value = $vals.R
}
